$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.251.59"
$ws.Range("E2").Value = "  -5.98%  "

$ws.Range("D3").Value = "2.858.82"
$ws.Range("E3").Value = "  -9.71%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Formula = "'461.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -12.88%  "

$ws.Range("D6").Formula = "'123.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.56%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "2.860.35"
$ws.Range("E8").Value = "  -9.57%  "

$ws.Range("E9").Value = "  -11.79%  "

$ws.Range("D10").Formula = "'6.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.07%  "

$ws.Range("D11").Formula = "'0.0943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -15.25%  "

$ws.Range("D12").Formula = "'0.321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -18.52%  "

$ws.Range("D13").Formula = "'0.121"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.25%  "

$ws.Range("D14").Value = "3.349.02"
$ws.Range("E14").Value = "  -9.87%  "

$ws.Range("D15").Formula = "'22.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -12.12%  "

$ws.Range("D16").Value = "54.250.53"
$ws.Range("E16").Value = "  -6.25%  "

$ws.Range("D17").Value = "2.855.83"
$ws.Range("E17").Value = "  -10.05%  "

$ws.Range("E18").Value = "  -15.84%  "

$ws.Range("E19").Value = "  -10.46%  "

$ws.Range("D20").Formula = "'11.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -15.75%  "

$ws.Range("D21").Formula = "'6.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -14.55%  "

$ws.Range("D22").Formula = "'290.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -18.74%  "

$ws.Range("D23").Formula = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Formula = "'0.431"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -16.50%  "

$ws.Range("D25").Formula = "'57.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -16.97%  "

$ws.Range("D26").Formula = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  -12.14%  "

$ws.Range("D29").Value = "0.0₃0779"
$ws.Range("E29").Value = "  -18.63%  "

$ws.Range("E30").Value = "  -14.00%  "

$ws.Range("E31").Value = "  -13.76%  "

$ws.Range("D32").Formula = "'1.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.18%  "

$ws.Range("D33").Formula = "'1.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.44%  "

$ws.Range("D34").Formula = "'18.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -15.89%  "

$ws.Range("D35").Formula = "'136.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -14.41%  "

$ws.Range("E36").Value = "  -18.24%  "

$ws.Range("D37").Formula = "'5.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -16.12%  "

$ws.Range("D38").Formula = "'1.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -17.20%  "

$ws.Range("D39").Formula = "'22.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.16%  "

$ws.Range("D40").Value = "2.885.93"
$ws.Range("E40").Value = "  -9.82%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").Formula = "'0.0604"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -14.37%  "

$ws.Range("E43").Value = "  -12.99%  "

$ws.Range("E44").Value = "  -15.11%  "

$ws.Range("D45").Formula = "'0.915"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -16.20%  "

$ws.Range("E46").Value = "  -13.40%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.015.91"
$ws.Range("E47").Value = "  -11.28%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Formula = "'3.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -17.56%  "

$ws.Range("D49").Formula = "'5.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -15.57%  "

$ws.Range("E50").Value = "  -11.53%  "

$ws.Range("D51").Formula = "'17.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -16.03%  "
